# AutoCommit_13 ноября 2023 г. 16:21:54_SibNout2023
# Fill in grades that were entered for several students (rows 7, 21, 22, 26, 27, 29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Арутюнян Алик
$ws.Range("E7").Value = 5

# Row 21 - Мохначев Егор
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 5

# Row 22 - Мукминова Диана
$ws.Range("E22").Value = 5

# Row 26 - Поздин Александр
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 4

# Row 27 - Райлян Арсений
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 5
$ws.Range("I27").Value = 5

# Row 29 - Тумат Ангыр
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 5

# Final selection left on I30 after entering the last grade
$ws.Range("I30").Select()
